$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-05 Tuesday" "2024-11-06 Wednesday"
Replace-Text "52×93=" "29×34="
Replace-Text "46×11=" "90×84="
Replace-Text "20×59=" "86×71="
Replace-Text "52×20=" "43×71="
Replace-Text "28×82=" "13×26="
Replace-Text "97×13=" "29×40="
Replace-Text "76×65=" "92×86="
Replace-Text "24×68=" "97×16="
Replace-Text "87×65=" "69×70="
Replace-Text "33×27=" "54×98="
Replace-Text "55×45=" "77×36="
Replace-Text "83×35=" "66×97="
Replace-Text "74×48=" "60×49="
Replace-Text "15×25=" "86×57="
Replace-Text "78×20=" "45×58="
Replace-Text "64×43=" "77×91="
Replace-Text "48×86=" "41×75="
Replace-Text "18×87=" "59×99="
Replace-Text "11×31=" "67×59="
Replace-Text "46×49=" "38×21="
Replace-Text "20×63=" "47×40="
Replace-Text "81×88=" "91×98="
Replace-Text "82×86=" "82×84="
Replace-Text "47×36=" "45×46="
Replace-Text "21×40=" "91×42="
